$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("DPLKKPS132-001")
$ws2 = $wb.Worksheets.Item("DPLKKPS132-002")

# Update the register number referenced on the "Kembalikan ke Register" sheet
# (DPLKKPS132-002) first, then the "Setuju" sheet (DPLKKPS132-001) -- this
# write order matches the shared-string table ordering seen in the target file.
$ws2.Range("N2").Value = "M13220800000023"
$ws2.Range("F2").Value = "Username : 31816;`nPassword : bni1234;`nRole : Penyelia Settlement;`nNo. Register : M13220800000023;`nStatus Verifikasi : 0 : Kembalikan ke Register;`nKeterangan Verifikasi : KEP.TRX.436 tidak disetujui"

$ws1.Range("N2").Value = "M13220800000023"
$ws1.Range("F2").Value = "Username : 31816;`nPassword : bni1234;`nRole : 09;`nNo. Register : M13220800000023;`nStatus Verifikasi : 1 : Setuju;`nKeterangan Verifikasi : KEP.TRX.436 Setuju"

# Move the active tab / selection back to sheet 1 (DPLKKPS132-001) and update
# each sheet's selected cell to G2.
$ws2.Range("G2").Select()
$ws1.Activate()
$ws1.Range("G2").Select()
